# The deck currently has the "Integral" theme applied to the slide master
# (ppt/theme/theme2.xml, the theme actually driving every slide/layout) while
# the original default "Office Theme" colours were left behind as the theme
# used only by the Notes Master (ppt/theme/theme1.xml). The edit swaps the
# two: the presentation's live colour theme reverts to the stock "Office"
# palette.
#
# PowerPoint's ThemeColorScheme exposes the 12 DrawingML theme colour slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) for the theme that backs
# the slide master/layouts/slides - i.e. ppt/theme/theme2.xml in this deck.
# Driving every slot back to the standard Office values reproduces the
# colour swap described by the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# index : scheme slot : target "Office" RGB
$officeColors = @{
    1  = 0          # dk1      #000000
    2  = 16777215   # lt1      #FFFFFF
    3  = 6968388     # dk2      #44546A
    4  = 15132391    # lt2      #E7E6E6
    5  = 13998939    # accent1  #5B9BD5
    6  = 3243501     # accent2  #ED7D31
    7  = 10855845    # accent3  #A5A5A5
    8  = 49407       # accent4  #FFC000
    9  = 12874308    # accent5  #4472C4
    10 = 4697456     # accent6  #70AD47
    11 = 12673797    # hlink    #0563C1
    12 = 7491477     # folHlink #954F72
}

foreach ($idx in 1..12) {
    $tcs.Colors($idx).RGB = $officeColors[$idx]
}
